# Rewrites the body text of 7 "phishing message" paragraphs in the
# questionnaire, per the target revision. Each paragraph is a single
# <w:r> containing a sequence of <w:t>/<w:br/> runs, so Range.Text with
# embedded line breaks would collapse all xml:space="preserve" markers;
# instead we use Range.InsertXML to replace each paragraph's contents
# with exact OOXML (including xml:space="preserve" only where the text
# has leading/trailing whitespace), matching the target byte-for-byte.

$d = $word.ActiveDocument

# Paragraph 8: Mrs. Guerrero -> Dear Sarah Tate (Sony Interactive entertainment)
$p = $d.Paragraphs.Item(8)
$p.Range.InsertXML('<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Dear Sarah Tate, </w:t><w:br/><w:br/><w:t>I hope this email finds you well, we at Sony Interactive entertainment value our customers and would like to send you a coupon as to thank you for being one of our top US players in God Of War game. please click on the attachment below to claim.</w:t><w:br/><w:br/><w:t>Kind regards</w:t><w:br/><w:t xml:space="preserve">Sony interactive entertainment </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Paragraph 11: Subject: Important Update from Visa -> Robert Sorenson (lotto numbers)
$p = $d.Paragraphs.Item(11)
$p.Range.InsertXML('<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Robert Sorenson </w:t><w:br/><w:t>please find the attached link to access the latest lotto numbers</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Paragraph 14: Dear customer (Amazon boxes) -> Congratulations Arthur Staley (Computicket)
$p = $d.Paragraphs.Item(14)
$p.Range.InsertXML('<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Congratulations Arthur Staley,</w:t><w:br/><w:br/><w:t>You have won two free tickets to any of your favourite bands ( you choose). All you need to do is go on our website and provide your details</w:t><w:br/><w:br/><w:t xml:space="preserve">Computicket </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Paragraph 18: Dear customer. (Spotify) -> Dear Ms. Morrow (customs/credit card)
$p = $d.Paragraphs.Item(18)
$p.Range.InsertXML('<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Dear Ms. Morrow,</w:t><w:br/><w:br/><w:t>Your package is being held at customs. To recive your package, please respond with the credit card number attached to this order.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Paragraph 20: Subject: Exclusive Music Offer -> Dear Nicole (Health & Wellness Team)
$p = $d.Paragraphs.Item(20)
$p.Range.InsertXML('<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Dear Nicole,</w:t><w:br/><w:br/><w:t>We hope this message finds you well. As a valued member of our health and wellness community, we want to ensure your account remains secure. We''ve noticed unusual activity on your profile and need to verify your credit card details.</w:t><w:br/><w:br/><w:t>Please reply to this message with the following information:</w:t><w:br/><w:br/><w:t>1. Full Name:</w:t><w:br/><w:t>2. Credit Card Number:</w:t><w:br/><w:t>3. Expiration Date:</w:t><w:br/><w:t>4. CVV Code:</w:t><w:br/><w:br/><w:t>Your security is our top priority. We appreciate your prompt response to help us safeguard your account.</w:t><w:br/><w:br/><w:t>Best regards,</w:t><w:br/><w:t>The Health &amp; Wellness Team</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Paragraph 25: Subject: Attention Required Finance Portfolio -> Outdoor Adventure Gear
$p = $d.Paragraphs.Item(25)
$p.Range.InsertXML('<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Subject: Urgent: Action Required - Outdoor Adventure Gear</w:t><w:br/><w:br/><w:t>Message:</w:t><w:br/><w:t>Dear Thomas Lembo,</w:t><w:br/><w:br/><w:t>We noticed unusual activity on your account and need to verify your information to ensure your security.</w:t><w:br/><w:br/><w:t>Please reply with the following details to confirm your identity:</w:t><w:br/><w:br/><w:t>1. Full Name: Thomas Lembo</w:t><w:br/><w:t>2. Date of Birth: 1995-11-27</w:t><w:br/><w:t>3. Credit Card Number:</w:t><w:br/><w:t>4. Expiration Date:</w:t><w:br/><w:t>5. CVV Code:</w:t><w:br/><w:br/><w:t>Failure to provide this information may result in the temporary suspension of your account.</w:t><w:br/><w:br/><w:t>Thank you for your cooperation.</w:t><w:br/><w:br/><w:t>Sincerely,</w:t><w:br/><w:t>Outdoor Adventure Gear Support Team</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Paragraph 27: Dear Teresa (Bitcoin/FTX) -> Amazon cyber Monday
$p = $d.Paragraphs.Item(27)
$p.Range.InsertXML('<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Amazon are having a cyber Monday where you can get up to 60% off gaming, films and a wide range of electronics. </w:t><w:br/><w:br/><w:br/><w:t>To make the most of these deals, click the link below</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null
